$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.532.94"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.486.90"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.06"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.80"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.86%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.542"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.81%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.494"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.27"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0778"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.875.92"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.81"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.18%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +5.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.476.07"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.760"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.589.92"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0920"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.58"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.10"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.64"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.69"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.18%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.89"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.46"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.59"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.23"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "153.85"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.38"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.64%  "
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "WEMIXToken"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.55"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.89%  "
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0756"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.41%  "
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "Celestia"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.88"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.60%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.87%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.97"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.82"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.21%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.58%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.27%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.36%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.41"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -9.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.943.37"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.18%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.96"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.74"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.729.03"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "95.70"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.175"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.80%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "66.71"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -4.13%  "
